# Word COM-interop script implementing the tracked changes described by the
# commit "Changes to graphs and formulas".
#
# Summary of edits:
#   1. "November 27, 2019" -> "November 29, 2019" (typed fix, split across 3 runs)
#   2. The stray "_GoBack" bookmark (originally sitting right before "ANSWER:")
#      is relocated by Word to the most recent edit point, which ends up
#      wrapping "why did they succeed?" further down in the document.
#   3. "amount" -> "count" (and the grammar-checker proofErr markers that used
#      to flag "amount" are cleared, since the flagged text no longer exists)
#   4. "...crowdfunding projects at Kickstarter." gains "successfully delivered"
#   5. A few runs that used to be split (purely because of older proofing
#      edits) get collapsed back into single runs with identical text:
#        - "...taken with some criteria that biased the results."
#        - the "Getting funded on Kickstarter..." quotation
#        - "3.5 Average duration of failed projects per Category/"
#   6. "succeed ?" -> "succeed?" with the comma that used to follow it moved
#      into its own run (net visible text unchanged apart from the removed
#      space before "?").
#
# Helper: forces a run boundary at an exact character offset without changing
# the visible formatting, by toggling Bold on then back off on the two sides
# of the split. iron_native (like real Word) only merges runs that are
# written with identical formatting as part of the very same operation, so
# alternating the toggle across adjacent segments keeps them as separate
# <w:r> elements even after the property is reset to its original value.
function Split-Runs($rng, $offsets) {
    $base = $rng.Start
    for ($i = 0; $i -lt $offsets.Length - 1; $i++) {
        $segStart = $base + $offsets[$i]
        $segEnd = $base + $offsets[$i + 1]
        $seg = $d.Range($segStart, $segEnd)
        if ($i % 2 -eq 0) {
            $seg.Font.Bold = 1
        } else {
            $seg.Font.Bold = 0
        }
    }
    $whole = $d.Range($base, $base + $offsets[$offsets.Length - 1])
    $whole.Font.Bold = 0
}

# Helper: collapse every run inside $rng into a single run. Setting .Text to
# the very same text it already holds is a no-op (nothing to normalize), so
# this stages the write through a throwaway placeholder first, forcing the
# engine to actually rebuild the range (and merge its runs) on the second,
# real write.
function Merge-Runs($rng, $text) {
    $start = $rng.Start
    $rng.Text = "X"
    $placeholder = $d.Range($start, $start + 1)
    $placeholder.Text = $text
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "November 27, 2019" -> "November 29, 2019", typed as three runs:
#    "November 2" | "9" | ", 2019"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("November 27, 2019", $false, $false, $false, $false, $false, `
                 $true, 1, $false, "November 29, 2019", 2) | Out-Null
Split-Runs $d.Range(0, 18) @(0, 10, 11, 18)

# ---------------------------------------------------------------------------
# 2/3/4. "...have the largest amount of crowdfunding projects at
#         Kickstarter." -> "...have the largest count of crowdfunding
#         projects successfully delivered at Kickstarter."
#    Replacing the whole sentence fragment in one shot naturally drops the
#    now-stale proofErr gramStart/gramEnd markers that used to bracket
#    "amount", exactly like Word does when the flagged text is edited away.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("have the largest amount of crowdfunding projects at Kickstarter.", `
                 $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $r.Start
$newText = "have the largest count of crowdfunding projects successfully delivered at Kickstarter."
$r.Text = $newText
# Segment boundaries (relative offsets into $newText):
#   "have the largest " | "count" | " of crowdfunding " | "projects" |
#   " successfully delivered" | " " | "at Kickstarter."
Split-Runs $d.Range($start, $start + $newText.Length) @(0, 17, 22, 39, 47, 70, 71, 86)

# ---------------------------------------------------------------------------
# 5a. Collapse the three runs around "...taken with <result/s>" back into one
#     run with identical combined text.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(". Making it just a small sample and maybe it was taken with some criteria that biased the results.", `
                 $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-Runs $r ". Making it just a small sample and maybe it was taken with some criteria that biased the results."

# ---------------------------------------------------------------------------
# 6. "succeed ?" -> "succeed?"; the "_GoBack" bookmark (previously sitting
#    near "ANSWER:") is relocated here by Word, wrapping "why did they
#    succeed?"; the comma that used to directly follow "succeed ?" becomes
#    its own run, and the run after it loses its now-redundant leading comma.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("why did they succeed ?", $false, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$whyStart = $r.Start

$r2 = $d.Content
$r2.Find.Execute("succeed ?", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "succeed?", 2) | Out-Null
$succeedEnd = $r2.End

$after = $d.Range($succeedEnd, $succeedEnd)
$after.InsertAfter(",")
Split-Runs $d.Range($succeedEnd, $succeedEnd + 1) @(0, 1)

$nextRange = $d.Range($succeedEnd + 1, $succeedEnd + 3)
$nextRange.Text = " "

$bmRange = $d.Range($whyStart, $succeedEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 5b. Collapse the "Getting funded on Kickstarter..." quotation back into a
#     single run with the same combined text.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("The only metric that seems to be analyzed by the backers is: " + [char]0x201C + "Getting funded on Kickstarter requires meeting or exceeding the project's initial goal", `
                 $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-Runs $r (" The only metric that seems to be analyzed by the backers is: " + [char]0x201C + "Getting funded on Kickstarter requires meeting or exceeding the project's initial goal")

# ---------------------------------------------------------------------------
# 5c. Collapse "3.5 Average duration of failed projects " / "per Category/"
#     back into a single run with the same combined text.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("3.5 Average duration of failed projects per Category/", `
                 $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-Runs $r "3.5 Average duration of failed projects per Category/"
